# Commit: "Store the user profile in project path (EmployeeId is a folder name)"
# Adds a new employee record (row 13) to the "Employees" sheet, fills in the
# previously-blank "Current Role" for the existing last row (row 12), and
# leaves the "Roles" sheet content logically unchanged (its shared string
# references merely shift because new unique strings were introduced earlier
# in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")

# Row 12 (E0123 / p23) previously had an empty "Current Role" cell -> now
# filled in with "Solution Developer".
$ws.Range("F12").Value = "Solution Developer"

# Duplicate row 12's formatting down into the brand new row 13 so that every
# cell keeps the same number formats / styles used throughout the table.
$ws.Range("A12:K12").Copy()
$ws.Range("A13:K13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(13).RowHeight = $ws.Rows.Item(12).RowHeight

# Populate the new employee row (Sr No 12 / EMP ID E0124 / p24).
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "E0124"
$ws.Range("C13").Value = "p24"
$ws.Range("D13").Value = 45287
$ws.Range("E13").Value = 36329
$ws.Range("F13").Value = "Solution Developer"
$ws.Range("G13").Value = "p24@gmail.com"
$ws.Range("H13").Value = "Female"
$ws.Range("I13").Value = 1234543268
$ws.Range("J13").Value = "Pune"
$ws.Range("K13").Value = 1

# Match the workbook's recorded selection after the edit.
$ws.Activate() | Out-Null
$ws.Range("G13").Select() | Out-Null
